$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 07:03"

# Apply updated country stats / re-sorted rows
# Row 21 -> Suiza
$ws.Cells.Item(21, 1).Value = "Suiza"
$ws.Cells.Item(21, 2).Value = 30060
$ws.Cells.Item(21, 4).Value = 25700
$ws.Cells.Item(21, 5).Value = 2555
$ws.Cells.Item(21, 6).Value = 121
$ws.Cells.Item(21, 8).Value = 1805

# Row 22 -> Ecuador
$ws.Cells.Item(22, 1).Value = "Ecuador"
$ws.Cells.Item(22, 2).Value = 29420
$ws.Cells.Item(22, 4).Value = 3433
$ws.Cells.Item(22, 5).Value = 24369
$ws.Cells.Item(22, 6).Value = 156
$ws.Cells.Item(22, 8).Value = 1618

# Row 64 -> Hungria
$ws.Cells.Item(64, 1).Value = "Hungria"
$ws.Cells.Item(64, 2).Value = 3150
$ws.Cells.Item(64, 3).Value = 39
$ws.Cells.Item(64, 4).Value = 801
$ws.Cells.Item(64, 5).Value = 1966
$ws.Cells.Item(64, 6).Value = 50
$ws.Cells.Item(64, 7).Value = 10
$ws.Cells.Item(64, 8).Value = 383

# Row 65 -> Nigeria
$ws.Cells.Item(65, 1).Value = "Nigeria"
$ws.Cells.Item(65, 2).Value = 3145
$ws.Cells.Item(65, 4).Value = 534
$ws.Cells.Item(65, 5).Value = 2508
$ws.Cells.Item(65, 6).Value = 4
$ws.Cells.Item(65, 8).Value = 103

# Row 67
$ws.Cells.Item(67, 2).Value = 2992
$ws.Cells.Item(67, 3).Value = 3
$ws.Cells.Item(67, 4).Value = 2772
$ws.Cells.Item(67, 5).Value = 165

# Row 95 -> Kirguistan
$ws.Cells.Item(95, 1).Value = "Kirguistan"
$ws.Cells.Item(95, 2).Value = 895
$ws.Cells.Item(95, 3).Value = 24
$ws.Cells.Item(95, 4).Value = 637
$ws.Cells.Item(95, 5).Value = 246
$ws.Cells.Item(95, 6).Value = 13
$ws.Cells.Item(95, 8).Value = 12

# Row 96 -> Republica de Chipre
$ws.Cells.Item(96, 1).Value = "Republica de Chipre"
$ws.Cells.Item(96, 2).Value = 883
$ws.Cells.Item(96, 4).Value = 296
$ws.Cells.Item(96, 5).Value = 572
$ws.Cells.Item(96, 6).Value = 15
$ws.Cells.Item(96, 8).Value = 15

# Row 97 -> Somalia
$ws.Cells.Item(97, 1).Value = "Somalia"
$ws.Cells.Item(97, 2).Value = 873
$ws.Cells.Item(97, 4).Value = 87
$ws.Cells.Item(97, 5).Value = 747
$ws.Cells.Item(97, 6).Value = 2
$ws.Cells.Item(97, 8).Value = 39

# Row 110 -> El Salvador
$ws.Cells.Item(110, 1).Value = "El Salvador"
$ws.Cells.Item(110, 2).Value = 695
$ws.Cells.Item(110, 3).Value = 62
$ws.Cells.Item(110, 4).Value = 245
$ws.Cells.Item(110, 5).Value = 435
$ws.Cells.Item(110, 6).Value = 4
$ws.Cells.Item(110, 8).Value = 15

# Row 111 -> Uruguay
$ws.Cells.Item(111, 1).Value = "Uruguay"
$ws.Cells.Item(111, 2).Value = 673
$ws.Cells.Item(111, 4).Value = 486
$ws.Cells.Item(111, 5).Value = 170
$ws.Cells.Item(111, 6).Value = 10
$ws.Cells.Item(111, 8).Value = 17

# Row 191 -> Nueva Caledonia
$ws.Cells.Item(191, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(191, 4).Value = 18
$ws.Cells.Item(191, 8).Value = 0

# Row 192 -> Belice
$ws.Cells.Item(192, 1).Value = "Belice"
$ws.Cells.Item(192, 4).Value = 16
$ws.Cells.Item(192, 8).Value = 2

# Row 198 -> Dominica
$ws.Cells.Item(198, 1).Value = "Dominica"
$ws.Cells.Item(198, 4).Value = 14
$ws.Cells.Item(198, 8).Value = 0

# Row 199 -> Curazao
$ws.Cells.Item(199, 1).Value = "Curazao"
$ws.Cells.Item(199, 4).Value = 13
$ws.Cells.Item(199, 8).Value = 1
